$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts existing B,C,D,E -> C,D,E,F
$ws.Columns.Item(2).Insert()

# Copy header formatting (bold/border/alignment style) from the neighboring header cell (now C1)
# onto the new header cell B1, then set its text.
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 2).Value = "segments"

# List of segment labels (previously in column A, now destined for column B)
$labels = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    # Column A becomes a numeric segment index starting at 0 (keeps existing bold/border style)
    $ws.Cells.Item($row, 1).Value = $i
    # Column B gets the text label with no special style (matching the diff, which has no s= attr)
    $ws.Cells.Item($row, 2).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Style = "Normal"
}

$ws.Range("A1").Select()
